$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.1243246666666667
$ws.Range("H2").Value = 0.372974
$ws.Range("I2").Value = 0.09963085929726231
$ws.Range("J2").Value = 0.09963085929726233
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.155977333333333
$ws.Range("N2").Value = 9.467931999999999
$ws.Range("O2").Value = 0.3579027849973545
$ws.Range("P2").Value = 0.3579027849973545
$ws.Range("Q2").Value = 0.3923658299742222
$ws.Range("R2").Value = 3.531292469768
$ws.Range("S2").Value = 0.03565816201416975
$ws.Range("T2").Value = 0.03565816201416975

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.1243246666666667
$ws.Range("H3").Value = 0.372974
$ws.Range("I3").Value = 0.09963085929726231
$ws.Range("J3").Value = 0.09963085929726233
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.165953666666667
$ws.Range("N3").Value = 9.497861
$ws.Range("O3").Value = 0.359034148472735
$ws.Range("P3").Value = 0.359034148472735
$ws.Range("Q3").Value = 0.3936061342904445
$ws.Range("R3").Value = 3.542455208614
$ws.Range("S3").Value = 0.03577088072939944
$ws.Range("T3").Value = 0.03577088072939945

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.1243246666666667
$ws.Range("H4").Value = 0.372974
$ws.Range("I4").Value = 0.09963085929726231
$ws.Range("J4").Value = 0.09963085929726233
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.496042666666666
$ws.Range("N4").Value = 7.488128
$ws.Range("O4").Value = 0.2830630665299106
$ws.Range("P4").Value = 0.2830630665299106
$ws.Range("Q4").Value = 0.3103196725191111
$ws.Range("R4").Value = 2.792877052672
$ws.Range("S4").Value = 0.02820181655369313
$ws.Range("T4").Value = 0.02820181655369313

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.7328223333333334
$ws.Range("H5").Value = 2.198467
$ws.Range("I5").Value = 0.5872665557027417
$ws.Range("J5").Value = 0.5872665557027417
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.155977333333333
$ws.Range("N5").Value = 9.467931999999999
$ws.Range("O5").Value = 0.3579027849973545
$ws.Range("P5").Value = 0.3579027849973545
$ws.Range("Q5").Value = 2.312770673360444
$ws.Range("R5").Value = 20.814936060244
$ws.Range("S5").Value = 0.2101843358218152
$ws.Range("T5").Value = 0.2101843358218152

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.7328223333333334
$ws.Range("H6").Value = 2.198467
$ws.Range("I6").Value = 0.5872665557027417
$ws.Range("J6").Value = 0.5872665557027417
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.165953666666667
$ws.Range("N6").Value = 9.497861
$ws.Range("O6").Value = 0.359034148472735
$ws.Range("P6").Value = 0.359034148472735
$ws.Range("Q6").Value = 2.320081553231889
$ws.Range("R6").Value = 20.880733979087
$ws.Range("S6").Value = 0.2108487477532498
$ws.Range("T6").Value = 0.2108487477532498

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.7328223333333334
$ws.Range("H7").Value = 2.198467
$ws.Range("I7").Value = 0.5872665557027417
$ws.Range("J7").Value = 0.5872665557027417
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.496042666666666
$ws.Range("N7").Value = 7.488128
$ws.Range("O7").Value = 0.2830630665299106
$ws.Range("P7").Value = 0.2830630665299106
$ws.Range("Q7").Value = 1.829155811086222
$ws.Range("R7").Value = 16.462402299776
$ws.Range("S7").Value = 0.1662334721276766
$ws.Range("T7").Value = 0.1662334721276766

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.2213483333333333
$ws.Range("H8").Value = 0.664045
$ws.Range("I8").Value = 0.1773833402919521
$ws.Range("J8").Value = 0.1773833402919521
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 3.155977333333333
$ws.Range("N8").Value = 9.467931999999999
$ws.Range("O8").Value = 0.3579027849973545
$ws.Range("P8").Value = 0.3579027849973545
$ws.Range("Q8").Value = 0.6985703227711111
$ws.Range("R8").Value = 6.287132904939999
$ws.Range("S8").Value = 0.0634859915026231
$ws.Range("T8").Value = 0.0634859915026231

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.2213483333333333
$ws.Range("H9").Value = 0.664045
$ws.Range("I9").Value = 0.1773833402919521
$ws.Range("J9").Value = 0.1773833402919521
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.165953666666667
$ws.Range("N9").Value = 9.497861
$ws.Range("O9").Value = 0.359034148472735
$ws.Range("P9").Value = 0.359034148472735
$ws.Range("Q9").Value = 0.7007785675272222
$ws.Range("R9").Value = 6.307007107745
$ws.Range("S9").Value = 0.06368667653497041
$ws.Range("T9").Value = 0.06368667653497041

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.2213483333333333
$ws.Range("H10").Value = 0.664045
$ws.Range("I10").Value = 0.1773833402919521
$ws.Range("J10").Value = 0.1773833402919521
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.496042666666666
$ws.Range("N10").Value = 7.488128
$ws.Range("O10").Value = 0.2830630665299106
$ws.Range("P10").Value = 0.2830630665299106
$ws.Range("Q10").Value = 0.5524948841955555
$ws.Range("R10").Value = 4.97245395776
$ws.Range("S10").Value = 0.05021067225435862
$ws.Range("T10").Value = 0.05021067225435862

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.1693576666666667
$ws.Range("H11").Value = 0.508073
$ws.Range("I11").Value = 0.1357192447080439
$ws.Range("J11").Value = 0.1357192447080439
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 3.155977333333333
$ws.Range("N11").Value = 9.467931999999999
$ws.Range("O11").Value = 0.3579027849973545
$ws.Range("P11").Value = 0.3579027849973545
$ws.Range("Q11").Value = 0.5344889572262221
$ws.Range("R11").Value = 4.810400615035999
$ws.Range("S11").Value = 0.04857429565874636
$ws.Range("T11").Value = 0.04857429565874636

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.1693576666666667
$ws.Range("H12").Value = 0.508073
$ws.Range("I12").Value = 0.1357192447080439
$ws.Range("J12").Value = 0.1357192447080439
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 3.165953666666667
$ws.Range("N12").Value = 9.497861
$ws.Range("O12").Value = 0.359034148472735
$ws.Range("P12").Value = 0.359034148472735
$ws.Range("Q12").Value = 0.5361785257614444
$ws.Range("R12").Value = 4.825606731853
$ws.Range("S12").Value = 0.04872784345511527
$ws.Range("T12").Value = 0.04872784345511527

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.1693576666666667
$ws.Range("H13").Value = 0.508073
$ws.Range("I13").Value = 0.1357192447080439
$ws.Range("J13").Value = 0.1357192447080439
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 2.496042666666666
$ws.Range("N13").Value = 7.488128
$ws.Range("O13").Value = 0.2830630665299106
$ws.Range("P13").Value = 0.2830630665299106
$ws.Range("Q13").Value = 0.422723961927111
$ws.Range("R13").Value = 3.804515657344
$ws.Range("S13").Value = 0.03841710559418225
$ws.Range("T13").Value = 0.03841710559418225
